$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.065.67'
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = '1.637.51'
$ws.Range("E3").Value = '  -1.75%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '213.87'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.21%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.5237'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("E7").Value = '  -0.09%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2594'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.72%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06285'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.06%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '20.60'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.30%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07668'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.97%  '

$ws.Range("D12").Value = '1.645.53'
$ws.Range("E12").Value = '  -1.28%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.399'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.63%  '

$ws.Range("D14").Value = '1.860.61'
$ws.Range("E14").Value = '  -1.77%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.5505'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").Value = '0.0₅8214'
$ws.Range("E16").Value = '  +3.71%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '64.83'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.28%  '

$ws.Range("D18").Value = '26.056.13'
$ws.Range("E18").Value = '  -0.23%  '

$ws.Range("E19").Value = '  -0.07%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.674'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.77%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '187.97'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.01%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '10.18'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.65%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.148'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.10%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '145.41'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.60%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.1208'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.47%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.393'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.71%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '15.78'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.58%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.393'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.94%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.05949'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -5.93%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.254'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.56%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.425'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.80%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.394'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.43%  '

$ws.Range("E34").Value = '  +0.41%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.9804'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.96%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.394'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.52%  '

$ws.Range("E37").Value = '  +1.07%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.5630'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -6.32%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01613'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.02%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.8498'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.33%  '

$ws.Range("E41").Value = '  -0.21%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.695'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -6.47%  '

$ws.Range("D43").Value = '1.029.09'
$ws.Range("E43").Value = '  -6.99%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '100.11'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.30%  '

$ws.Range("D45").Value = '1.786.95'

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '55.66'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.74%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.9984'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.36%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '8.060'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.43%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.05143'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.69%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.4217'
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '5.896'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.46%  '
